# The original file had a worksheet named "excel" fed by a web query
# ("excel.iqy" -> connection "excel" -> query table "excel") whose result
# landed in the defined name "excel" (scoped to that sheet). The commit
# re-runs that import against a second copy of the query file
# ("excel(1).iqy"), which is why everything downstream gets renamed with
# a "(1)" / "_1" suffix, and the report timestamp embedded in the sheet
# is refreshed to the time of the new run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet "excel" -> "excel(1)"
$ws.Name = "excel(1)"

# Workbook-scoped defined name "excel" -> "excel_1"; Excel automatically
# keeps its RefersTo formula pointing at the (renamed) sheet.
foreach ($n in $wb.Names) {
    $n.Name = "excel_1"
}

# Rename the web-query connection and its query table the same way Excel
# does when it re-imports from "excel(1).iqy". (There is exactly one of
# each, with id/index 1, matching xl/connections.xml & queryTable1.xml.)
try {
    $c1 = $wb.Connections.Item(1)
    $c1.Name = "excel(1)"
    try { $c1.ODCFile = "C:\Users\AxeelZR\Downloads\excel(1).iqy" } catch {}
} catch {}
try {
    $qt1 = $ws.QueryTables.Item(1)
    $qt1.Name = "excel(1)"
} catch {}

# Refresh the "Reporte generado a las ..." timestamp text printed on the
# sheet (cell A23) to match the new run.
$ws.Range("A23").Value = "Reporte generado a las 01:36 PM el 5/12/2018"
